$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for all existing data rows
#    (rows 2 through 529) from 45189 (2023-09-20) to 45190 (2023-09-21).
$ws.Range("C2:C529").Value = 45190

# 2. Row 529 gains an explicit (default) row height, matching the target file.
$ws.Rows.Item(529).RowHeight = 15

# 3. Append two new data rows (530 and 531) for newly reported cases.

# --- Row 530 ---
$ws.Cells.Item(530, 1).Value = "A 44553-2023"

$ws.Cells.Item(530, 2).Value = 45189
$ws.Cells.Item(530, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(530, 3).Value = 45190
$ws.Cells.Item(530, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(530, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(530, 5).Value = "KRISTIANSTAD"
$ws.Cells.Item(530, 7).Value = 0.7
$ws.Cells.Item(530, 8).Value = 0
$ws.Cells.Item(530, 9).Value = 0
$ws.Cells.Item(530, 10).Value = 0
$ws.Cells.Item(530, 11).Value = 0
$ws.Cells.Item(530, 12).Value = 0
$ws.Cells.Item(530, 13).Value = 0
$ws.Cells.Item(530, 14).Value = 0
$ws.Cells.Item(530, 15).Value = 0
$ws.Cells.Item(530, 16).Value = 0
$ws.Cells.Item(530, 17).Value = 0
$ws.Cells.Item(530, 18).Value = ""
$ws.Cells.Item(530, 18).WrapText = $true

$ws.Rows.Item(530).RowHeight = 15

# --- Row 531 ---
$ws.Cells.Item(531, 1).Value = "A 44550-2023"

$ws.Cells.Item(531, 2).Value = 45189
$ws.Cells.Item(531, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(531, 3).Value = 45190
$ws.Cells.Item(531, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(531, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(531, 5).Value = "KRISTIANSTAD"
$ws.Cells.Item(531, 7).Value = 2
$ws.Cells.Item(531, 8).Value = 0
$ws.Cells.Item(531, 9).Value = 0
$ws.Cells.Item(531, 10).Value = 0
$ws.Cells.Item(531, 11).Value = 0
$ws.Cells.Item(531, 12).Value = 0
$ws.Cells.Item(531, 13).Value = 0
$ws.Cells.Item(531, 14).Value = 0
$ws.Cells.Item(531, 15).Value = 0
$ws.Cells.Item(531, 16).Value = 0
$ws.Cells.Item(531, 17).Value = 0
$ws.Cells.Item(531, 18).Value = ""
$ws.Cells.Item(531, 18).WrapText = $true

Write-Output "done"
